$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before AI (so old AI becomes AJ, shifting things right),
# which makes room for the new "MgCa Coretop modelled temperature" column.
$ws.Range("AI:AI").Insert()

# New header for inserted AI column
$ws.Range("AI1").Value = "MgCa Coretop modelled temperature"

# New header for AK (beyond the old last column which is now AJ)
$ws.Range("AK1").Value = "MgCa Temperature anomaly_BAYMAG - Coretop"

# Row 2 updated values
$ws.Range("W2").Value = 28.44
$ws.Range("X2").Value = 0.4449956936306414
$ws.Range("Y2").Value = 2.11188818363064
$ws.Range("Z2").Value = -0.7341709730360613
$ws.Range("AA2").Value = -1.308487639702662

$ws.Range("AI2").Value = 27.6266
$ws.Range("AJ2").Value = 0.07923333000000099
$ws.Range("AK2").Value = -0.4950833299999999
